$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two additional NXdata rows for dark-current / open-beam ("flat") images,
# mirroring the existing "data" row (row 5). Values are written in the same
# order the original author entered them so new shared-string entries line
# up with the canonical file.
$ws.Range("B6").Value = "data_dark"
$ws.Range("B7").Value = "data_flat"
$ws.Range("D6").Value = "image_dark"
$ws.Range("D7").Value = "image_flat"

$ws.Range("A6").Value = "NXdata"
$ws.Range("C6").Value = "odin_topic"
$ws.Range("E6").Value = "ADAr"
$ws.Range("F6").Value = "uint32"
$ws.Range("H6").Value = "480, 290, 3"

$ws.Range("A7").Value = "NXdata"
$ws.Range("C7").Value = "odin_topic"
$ws.Range("E7").Value = "ADAr"
$ws.Range("F7").Value = "uint32"
$ws.Range("H7").Value = "480, 290, 3"

# New column I: "custom_field" flag for each NX group.
$ws.Range("I6").Value = "yes"
$ws.Range("I7").Value = "yes"
$ws.Range("I2").Value = "no"
$ws.Range("I3").Value = "no"
$ws.Range("I4").Value = "no"
$ws.Range("I5").Value = "no"
$ws.Range("I1").Value = "custom_field"

# Header style for the new column: red, centered (like the other headers).
$ws.Range("I1").Font.Color = 255
$ws.Range("I1").Font.Name = "Calibri (Body)"
$ws.Range("I1").HorizontalAlignment = -4108

# Body cells use the same centered style as the rest of the table.
$ws.Range("H6:H7").HorizontalAlignment = -4108
$ws.Range("I2:I7").HorizontalAlignment = -4108

# Move the selection to I1, matching the saved workbook state.
$ws.Range("I1").Select()
